$d = $word.ActiveDocument

# Replace paragraph 6
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertXML(@"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Importância dos ensaios mecânicos; Introdução à extensometria; Medidas de dureza; Ensaio de tração; Ensaio de Flexão; Ensaios de compressão, Ensaios de torção, e resistência à tração por compressão diametral; Fratura dúctil e frágil e o ensaio de impacto Charpy e Izod.</w:t></w:r></w:p>
"@)

# Replace paragraph 7
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertXML(@"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>Significance of mechanical testing; Introduction to strain gage fundamentals; Hardness testing; Tensile testing; Bending testing; Standard methods for compression, torsion and splitting tensile strength determination; Ductile and Brittle fracture and the charpy and izod impact testing methods.</w:t></w:r></w:p>
"@)

# Replace paragraph 9
$p9 = $d.Paragraphs.Item(9)
$p9.Range.InsertXML(@"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>Descrever os ensaios mecânicos usuais para determinar propriedades de rigidez elástica, resistência e ductilidade dos</w:t><w:br/><w:t>materiais metálicos, cerâmicos e poliméricos. Descrever os principais métodos experimentais e técnicas de análise</w:t><w:br/><w:t>envolvidas. Habilitar os alunos ao tratamento e análise dos resultados experimentais obtidos pois a disciplina possui cunho teórico/prático com idas aos laboratórios para sedimentar a teoria vista em aulas expositivas.</w:t><w:br/></w:r><w:r><w:t>Programa</w:t><w:br/><w:br/><w:t>1.Importância dos ensaios e do estabelecimento de normas para sua execução; classificação dos ensaios. 2.Introdução à extensometria: análise experimental de tensões e deformações. 3.Medidas de dureza: Conceitos e tipos de ensaio. Ensaio de dureza por penetração: principais escalas e equipamentos. Ensaio de microdureza. Dureza ao choque. 4.O ensaio de tração. Parâmetros de ensaio. Propriedades obtidas no ensaio. Curvas tensão-deformação de engenharia e real. 5. Ensaios de flexão a três e quatro pontos: determinação do módulo de ruptura e módulo de Young, Fontes de erros no ensaio de flexão e como minimizá-los. Estatística de Weibull; 6. Ensaio de compressão: técnica experimental e cuidados necessários. Compressão de materiais dúcteis e frágeis. Resistência à tração por compressão diametral. Ensaio de torção: Determinação do módulo de elasticidade transversal. 7.Ensaio de impacto. Fatores que concorrem para a fratura frágil dos materiais. A transição dúctil-frágil. Métodos Charpy e Izod. Ensaio de impacto instrumentado. 8. Métodos experimentais para o ensaio de fluência. Formas de representação dos resultados. Taxa de fluência estacionária e determinação da energia de ativação para fluência.</w:t><w:br/></w:r><w:r><w:t>Duas avaliações escritas, compostas por provas e que poderão ser complementadas por trabalhos ou relatórios de experimentos realizados em laboratório.</w:t><w:br/></w:r><w:r><w:t>A cada avaliação (compreendendo uma prova, complementada por trabalho ou relatório) será atribuído grau entre zero e dez.</w:t></w:r></w:p>
"@)

# Replace paragraph 11
$p11 = $d.Paragraphs.Item(11)
$p11.Range.InsertXML(@"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Avaliação escrita. Para aprovação, a média entre a avaliação de Recuperação e o grau obtido no semestre deve ser maior ou igual a cinco.</w:t></w:r></w:p>
"@)

# Replace paragraph 12
$p12 = $d.Paragraphs.Item(12)
$p12.Range.InsertXML(@"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>The discipline is aimed at describing the usual mechanical testing methods to determine the properties of elastic stiffness, strength and ductility of metallic, ceramic and polymeric materials. Enable students to process and analyze the experimental results obtained, as the discipline has a theoretical/experimental nature with visits to laboratories to solidify the theory learned in lectures.</w:t></w:r></w:p>
"@)

# Replace paragraph 14
$p14 = $d.Paragraphs.Item(14)
$p14.Range.InsertXML(@"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>1.S.A. Souza. Ensaios Mecânicos de Materiais Metálicos. São Paulo: Editora Blucher. 5ª ed., 1982, 286p.</w:t><w:br/><w:t>2.A. Garcia, J.A. Spim, C.A. dos Santos. Ensaios dos Materiais. Rio de Janeiro: LTC Editora, 2ª ed., 2012, 384p.</w:t><w:br/><w:t>3.C.A. Sciammarella, F.M. Sciammarella. Mecânica Experimental dos Sólidos. Rio de Janeiro : LTC, 2017, 460p.</w:t><w:br/><w:t>4.R.W. Hertzberg. Deformation and Fracture Mechanics of Engineering Materials. New York: John Wiley &amp; Sons, 4th ed., 1996, 786p.</w:t><w:br/><w:t>5.C. Suryanarayana. Experimental Techniques in Materials and Mechanics. Boca Raton: CRC Press, 2011, 450p.</w:t><w:br/><w:t>6.N.E. Dowling, S.L. Kampe, M.V. Kral. Mechanical Behavior of Materials. Hoboken: Pearson, 5th ed., 2018, 946p.</w:t><w:br/><w:t>7.Y. Lee, J. Pan, R. Hathaway, M. Barkey. Fatigue Testing and Analysis Theory and Practice. Oxford: Elsevier Butterworth-Heinemann, 2005, 402p.</w:t><w:br/><w:t>8.R.W. Evans, B. Wilshire. Introduction to Creep. London: The Institute of Materials, 1993, 115p.</w:t><w:br/><w:t>9.L.F.M. Silva. Comportamento Mecânico dos Materiais. Porto: Publindústria Edições Técnicas, 2012, 391p.</w:t><w:br/><w:t>10.G.E. Dieter. Metalurgia Mecânica. Rio de Janeiro: Guanabara Dois, 2ª ed., 1981, 653p.</w:t></w:r></w:p>
"@)

# Replace paragraph 17
$p17 = $d.Paragraphs.Item(17)
$p17.Range.InsertXML(@"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Método: </w:t></w:r><w:r><w:t>471420 - Carlos Antonio Reis Pereira Baptista</w:t><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Critério: </w:t></w:r><w:r><w:t>3586455 - Cassius Olivio Figueiredo Terra Ruchert</w:t><w:br/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Norma de recuperação: </w:t></w:r><w:r><w:t>7459752 - Maria Ismenia Sodero Toledo Faria</w:t></w:r></w:p>
"@)

# Replace paragraph 19
$p19 = $d.Paragraphs.Item(19)
$p19.Range.InsertXML(@"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>5840793 - Sérgio Schneider</w:t></w:r></w:p>
"@)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
